$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A9").NumberFormat = "@"
$ws.Range("A9").Value = "01-09-2021"
$ws.Range("F9").Value = 3571500
